# books.xlsx update:
#  - append 3 new book rows (5-7)
#  - move the active selection to F5
#  - drop workbook protection element
#  - (sheet dimension / row count updates automatically from the data write)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of book data appended below the existing table.
$ws.Range("A5").Value = "Conclave"
$ws.Range("B5").Value = "Robert Harris"
$ws.Range("C5").Value = "Fiction"

$ws.Range("A6").Value = "Goodnight, Mister Tom"
$ws.Range("B6").Value = "Michelle Magorian"
$ws.Range("C6").Value = "Fiction"

$ws.Range("A7").Value = "Davita's Harp"
$ws.Range("B7").Value = "Chaim Potok"
$ws.Range("C7").Value = "Fiction"

# Workbook was unprotected (workbookProtection removed from workbook.xml).
[void]$wb.Unprotect()

# Active cell / selection moved to F5.
[void]$ws.Range("F5").Select()
